$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 currently holds Spain's data; we are inserting a new "Slovakia" row
# above it (alphabetically between Romania and Spain), which pushes Spain
# down to row 22 and drops the previously-empty row 22.

# 1) Copy the date formatting (column B/C style) from row 21 down to row 22
#    so the relocated Spain row keeps the left-aligned date format.
$ws.Range("B21:C21").Copy() | Out-Null
$ws.Range("B22:C22").PasteSpecial(-4122) | Out-Null

# 2) Move Spain's existing values down into row 22.
$ws.Range("A22").Value2 = $ws.Range("A21").Value2
$ws.Range("B22").Value2 = $ws.Range("B21").Value2
$ws.Range("C22").Value2 = $ws.Range("C21").Value2

# 3) Write the new Slovakia row into row 21.
$ws.Range("A21").Value2 = "Slovakia"
$ws.Range("B21").Value2 = [datetime]"2021-11-25"
$ws.Range("C21").Value2 = [datetime]"2020-12-09"

# 4) Match the saved selection/active cell.
$ws.Range("D22").Select() | Out-Null
